# "Work on "Introduction: About eFMI"."
#
# 1) The auto "last saved" date field (type="datetimeFigureOut") cached on the
#    Slide Master and on every Slide Layout moves from 16.05.2022 to
#    17.05.2022.
# 2) On slide 9 ("Legal information" / "History" overview) the two top-level
#    grouped boxes ("Gruppieren 1" = Legal information, "Gruppieren 29" =
#    History) get repositioned.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date-placeholder text everywhere it appears (master +
#    every custom layout), without touching any other placeholder.
# ---------------------------------------------------------------------------

function Update-DatePlaceholder($container, $newText) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master "17.05.2022"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout "17.05.2022"
}

# ---------------------------------------------------------------------------
# 2) Move the two grouped shapes on slide 9.
# ---------------------------------------------------------------------------

$slide9 = $p.Slides.Item(9)

for ($i = 1; $i -le $slide9.Shapes.Count; $i++) {
    $shp = $slide9.Shapes.Item($i)
    if ($shp.Name -eq "Gruppieren 1") {
        $shp.Left = 68.29070866141733
        $shp.Top = 150.16078950157467
    } elseif ($shp.Name -eq "Gruppieren 29") {
        $shp.Left = 384.87669381338605
        $shp.Top = 150.16078950157467
    }
}
